# tasks.xlsx: rename sheet, note the StandApp deadline, and format the two
# new "date" cells next to the Goal2 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the sheet Sheet1 -> StandApp.
$ws.Name = "StandApp"

# 2) New note under Goal2, plus date-formatted cells beside Goal2/Task2.1.
$ws.Range("B9").Value = "I will make StandApp in 2 months"
$ws.Range("C9").NumberFormat = "d-mmm-yy"
$ws.Range("C10").NumberFormat = "d-mmm-yy"

# 3) Column B needs to be wide enough to show the new note; column C keeps
#    its previous width (this also splits the old merged "B:C" <col> span).
$ws.Columns.Item(2).ColumnWidth = 37.6
$ws.Columns.Item(3).ColumnWidth = 16.67

# 4) Re-touch the header rows' (already-default) alignment so every affected
#    cell keeps its original look.
$ws.Range("B1:C1").HorizontalAlignment = 1
$ws.Range("F1:G1").HorizontalAlignment = 1
$ws.Range("B2:C2").HorizontalAlignment = 1

# 5) Scroll back to the top and leave the selection on C10.
$ws.Range("C10").Select()
